$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44483
$ws.Range("J2").Value = 120
$ws.Range("D3").Value = 44742
$ws.Range("H3").Value = 'Madrigal'
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 19500
$ws.Range("P3").Value = 488
$ws.Range("D4").Value = 44785
$ws.Range("H4").Value = 'Argentina(o)'
$ws.Range("I4").Value = 'Segunda'
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("N4").Value = '$/caja 50 unidades'
$ws.Range("P4").Value = 310
$ws.Range("Q4").Value = 50
$ws.Range("D5").Value = 44482
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("P5").Value = 362
$ws.Range("D6").Value = 44405
$ws.Range("K6").Value = 21000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21500
$ws.Range("P6").Value = 538
$ws.Range("D7").Value = 44468
$ws.Range("H7").Value = 'Argentina(o)'
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("N7").Value = '$/caja 50 unidades'
$ws.Range("P7").Value = 350
$ws.Range("Q7").Value = 50
$ws.Range("D8").Value = 44426
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19500
$ws.Range("P8").Value = 488
$ws.Range("D9").Value = 44762
$ws.Range("H9").Value = 'Madrigal'
$ws.Range("J9").Value = 200
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 488
$ws.Range("Q9").Value = 40
$ws.Range("D10").Value = 44827
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("P10").Value = 362
$ws.Range("D11").Value = 44356
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 19000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19500
$ws.Range("P11").Value = 390
$ws.Range("D12").Value = 44391
$ws.Range("J12").Value = 140
$ws.Range("D13").Value = 44160
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("P13").Value = 362
$ws.Range("Q13").Value = 40
$ws.Range("D14").Value = 44419
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 21000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 21500
$ws.Range("N14").Value = '$/caja 50 unidades'
$ws.Range("P14").Value = 430
$ws.Range("Q14").Value = 50
$ws.Range("D15").Value = 44363
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("O15").Value = 'Región de Coquimbo'
$ws.Range("P15").Value = 488
$ws.Range("Q15").Value = 40
$ws.Range("D16").Value = 44435
$ws.Range("H16").Value = 'Madrigal'
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 19000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 19500
$ws.Range("O16").Value = 'Región de Coquimbo'
$ws.Range("P16").Value = 488
$ws.Range("D17").Value = 44412
$ws.Range("H17").Value = 'Symphony'
$ws.Range("J17").Value = 240
$ws.Range("K17").Value = 21000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21500
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("P17").Value = 538
$ws.Range("Q17").Value = 40
$ws.Range("D18").Value = 44806
$ws.Range("H18").Value = 'Argentina(o)'
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14500
$ws.Range("O18").Value = 'Provincia de Limarí'
$ws.Range("P18").Value = 362
$ws.Range("D19").Value = 44370
$ws.Range("H19").Value = 'Argentina(o)'
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 21000
$ws.Range("M19").Value = 20429
$ws.Range("N19").Value = '$/caja 50 unidades'
$ws.Range("P19").Value = 409
$ws.Range("Q19").Value = 50
$ws.Range("D20").Value = 44370
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 22000
$ws.Range("L20").Value = 23000
$ws.Range("M20").Value = 22500
$ws.Range("P20").Value = 562
$ws.Range("D21").Value = 44377
$ws.Range("H21").Value = 'Madrigal'
$ws.Range("J21").Value = 150
$ws.Range("M21").Value = 20333
$ws.Range("N21").Value = '$/caja 40 unidades'
$ws.Range("P21").Value = 508
$ws.Range("Q21").Value = 40
$ws.Range("D22").Value = 44377
$ws.Range("H22").Value = 'Symphony'
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 21000
$ws.Range("L22").Value = 22000
$ws.Range("M22").Value = 21500
$ws.Range("P22").Value = 538
$ws.Range("D23").Value = 44706
$ws.Range("H23").Value = 'Madrigal'
$ws.Range("J23").Value = 250
$ws.Range("N23").Value = '$/caja 40 unidades'
$ws.Range("P23").Value = 538
$ws.Range("Q23").Value = 40
$ws.Range("D24").Value = 44769
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 17000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 17500
$ws.Range("P24").Value = 438
$ws.Range("D25").Value = 44384
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 21000
$ws.Range("L25").Value = 22000
$ws.Range("M25").Value = 21500
$ws.Range("P25").Value = 538
$ws.Range("D26").Value = 44384
$ws.Range("I26").Value = 'Segunda'
$ws.Range("J26").Value = 30
$ws.Range("M26").Value = 19333
$ws.Range("N26").Value = '$/caja 50 unidades'
$ws.Range("P26").Value = 387
$ws.Range("Q26").Value = 50
$ws.Range("D27").Value = 44384
$ws.Range("H27").Value = 'Symphony'
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 21000
$ws.Range("M27").Value = 20400
$ws.Range("P27").Value = 510
$ws.Range("D28").Value = 44398
$ws.Range("J28").Value = 170
$ws.Range("D29").Value = 44167
$ws.Range("H29").Value = 'Española'
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 13000
$ws.Range("L29").Value = 14000
$ws.Range("M29").Value = 13500
$ws.Range("N29").Value = '$/caja 30 unidades'
$ws.Range("O29").Value = 'Región Metropolitana'
$ws.Range("P29").Value = 450
$ws.Range("Q29").Value = 30
$ws.Range("D30").Value = 44433
$ws.Range("J30").Value = 160
$ws.Range("D31").Value = 44489
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 13000
$ws.Range("L31").Value = 14000
$ws.Range("M31").Value = 13500
